$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 18:15"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 3859296
$ws.Range("C4").Value = 26025
$ws.Range("D4").Value = 1776182
$ws.Range("E4").Value = 1940079
$ws.Range("G4").Value = 158
$ws.Range("H4").Value = 143035

# Row 6: 'India' -> 'India'
$ws.Range("B6").Value = 1113400
$ws.Range("C6").Value = 35536
$ws.Range("D6").Value = 696073
$ws.Range("E6").Value = 389855
$ws.Range("G6").Value = 644
$ws.Range("H6").Value = 27472

# Row 13: 'Reino Unido' -> 'Reino Unido'
$ws.Range("B13").Value = 294792
$ws.Range("C13").Value = 726
$ws.Range("G13").Value = 27
$ws.Range("H13").Value = 45300

# Row 17: 'Italia' -> 'Italia'
$ws.Range("B17").Value = 244434
$ws.Range("C17").Value = 218
$ws.Range("D17").Value = 196949
$ws.Range("E17").Value = 12440
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 35045

# Row 20: 'Alemania' -> 'Alemania'
$ws.Range("B20").Value = 202747
$ws.Range("C20").Value = 175
$ws.Range("E20").Value = 5785

# Row 24: 'Canada' -> 'Canada'
$ws.Range("B24").Value = 110329
$ws.Range("C24").Value = 330
$ws.Range("D24").Value = 97025
$ws.Range("E24").Value = 4452
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 8852

# Row 26: 'Irak' -> 'Irak'
$ws.Range("B26").Value = 92530
$ws.Range("C26").Value = 2310
$ws.Range("D26").Value = 60528
$ws.Range("E26").Value = 28221
$ws.Range("G26").Value = 90
$ws.Range("H26").Value = 3781

# Row 46: 'Singapur' -> 'Singapur'
$ws.Range("D46").Value = 44086
$ws.Range("E46").Value = 3799

# Row 51: 'Barein' -> 'Barein'
$ws.Range("E51").Value = 4114
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 125

# Row 56: 'Azerbaiyan' -> 'Azerbaiyan'
$ws.Range("B56").Value = 27521
$ws.Range("C56").Value = 388
$ws.Range("D56").Value = 18967
$ws.Range("E56").Value = 8200
$ws.Range("G56").Value = 5
$ws.Range("H56").Value = 354

# Row 61: 'Argelia' -> 'Argelia'
$ws.Range("B61").Value = 23084
$ws.Range("C61").Value = 535
$ws.Range("D61").Value = 16051
$ws.Range("E61").Value = 5955
$ws.Range("G61").Value = 10
$ws.Range("H61").Value = 1078

# Row 62: 'Serbia' -> 'Moldavia'
$ws.Range("A62").Value = "Moldavia"
$ws.Range("B62").Value = 20980
$ws.Range("C62").Value = 186
$ws.Range("D62").Value = 14376
$ws.Range("E62").Value = 5920
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 684

# Row 63: 'Moldavia' -> 'Serbia'
$ws.Range("A63").Value = "Serbia"
$ws.Range("B63").Value = 20894
$ws.Range("C63").Value = 396
$ws.Range("D63").Value = 14047
$ws.Range("E63").Value = 6375
$ws.Range("G63").Value = 11
$ws.Range("H63").Value = 472

# Row 70: 'Chequia' -> 'Chequia'
$ws.Range("B70").Value = 13902
$ws.Range("C70").Value = 47
$ws.Range("D70").Value = 8761
$ws.Range("E70").Value = 4782
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 359

# Row 72: 'Dinamarca' -> 'Kenia'
$ws.Range("A72").Value = "Kenia"
$ws.Range("B72").Value = 13353
$ws.Range("C72").Value = 603
$ws.Range("D72").Value = 4440
$ws.Range("E72").Value = 8688
$ws.Range("H72").Value = 225

# Row 73: 'Kenia' -> 'Dinamarca'
$ws.Range("A73").Value = "Dinamarca"
$ws.Range("B73").Value = 13173
$ws.Range("D73").Value = 12209
$ws.Range("E73").Value = 353
$ws.Range("H73").Value = 611

# Row 79: 'Etiopia' -> 'Republica de Macedonia'
$ws.Range("A79").Value = "Republica de Macedonia"
$ws.Range("B79").Value = 9153
$ws.Range("C79").Value = 127
$ws.Range("D79").Value = 4810
$ws.Range("E79").Value = 3921
$ws.Range("G79").Value = 8
$ws.Range("H79").Value = 422

# Row 80: 'Noruega' -> 'Etiopia'
$ws.Range("A80").Value = "Etiopia"
$ws.Range("B80").Value = 9147
$ws.Range("D80").Value = 2430
$ws.Range("E80").Value = 6554
$ws.Range("H80").Value = 163

# Row 81: 'Republica de Macedonia' -> 'Noruega'
$ws.Range("A81").Value = "Noruega"
$ws.Range("B81").Value = 9028
$ws.Range("D81").Value = 8138
$ws.Range("E81").Value = 635
$ws.Range("H81").Value = 255

# Row 85: 'Estado de Palestina' -> 'Estado de Palestina'
$ws.Range("E85").Value = 6567
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 61

# Row 96: 'Luxemburgo' -> 'Luxemburgo'
$ws.Range("B96").Value = 5605
$ws.Range("C96").Value = 122
$ws.Range("E96").Value = 1161

# Row 97: 'Republica de Yibuti' -> 'Republica de Yibuti'
$ws.Range("B97").Value = 5011
$ws.Range("C97").Value = 8
$ws.Range("D97").Value = 4838
$ws.Range("E97").Value = 117

# Row 101: 'Albania' -> 'Albania'
$ws.Range("B101").Value = 4090
$ws.Range("C101").Value = 82
$ws.Range("D101").Value = 2311
$ws.Range("E101").Value = 1667
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 112

# Row 102: 'Grecia' -> 'Grecia'
$ws.Range("B102").Value = 4007
$ws.Range("C102").Value = 24
$ws.Range("E102").Value = 2439

# Row 115: 'Mali' -> 'Mali'
$ws.Range("B115").Value = 2475
$ws.Range("C115").Value = 3
$ws.Range("D115").Value = 1828
$ws.Range("E115").Value = 526

# Row 118: 'Montenegro' -> 'Montenegro'
$ws.Range("B118").Value = 2188
$ws.Range("C118").Value = 116
$ws.Range("D118").Value = 385
$ws.Range("E118").Value = 1771
$ws.Range("G118").Value = 2
$ws.Range("H118").Value = 32

# Row 119: 'Estonia' -> 'Cabo Verde'
$ws.Range("A119").Value = "Cabo Verde"
$ws.Range("B119").Value = 2045
$ws.Range("C119").Value = 31
$ws.Range("D119").Value = 973
$ws.Range("E119").Value = 1051
$ws.Range("H119").Value = 21

# Row 120: 'Cabo Verde' -> 'Estonia'
$ws.Range("A120").Value = "Estonia"
$ws.Range("B120").Value = 2021
$ws.Range("D120").Value = 1912
$ws.Range("E120").Value = 40
$ws.Range("H120").Value = 69

# Row 137: 'Jordania' -> 'Jordania'
$ws.Range("B137").Value = 1218
$ws.Range("C137").Value = 4
$ws.Range("D137").Value = 1024
$ws.Range("E137").Value = 183

# Row 156: 'Malta' -> 'Malta'
$ws.Range("B156").Value = 677
$ws.Range("C156").Value = 1
$ws.Range("E156").Value = 6
